$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update card effect text (column D, rows 2-10).
# Rows 2-7 and 10 get their activation-cost wording reworded; row 8 gets a
# "消耗" prefix; row 9 (previously the "多选" card at row 9) swaps places
# content-wise with what is now row 8/10 per the authoring diff, but since
# we're editing by absolute cell, we just set the final text per row.
$ws.Range("D2").Value = "使用1任意属性，或将1张手牌洗回主牌堆：重抽房间区所有牌。使用1张《感知》发动本牌时，可以再重抽任意张手牌。"
$ws.Range("D3").Value = "使用1任意属性，或将1张手牌洗回主牌堆：选1张位于房间区最前方的怪物牌横置。使用1张《敏捷》发动本牌时，可以再选1张位于房间区最前方的怪物牌横置。"
$ws.Range("D4").Value = "使用1任意属性，或将1张手牌洗回主牌堆：将房间区最前方任意1张牌移动到战场敌人列第一行。使用1张《敏捷》发动本牌时，可以将房间区任意1张牌移动到战场敌人列第一行。"
$ws.Range("D5").Value = "使用1任意属性，或将1张手牌洗回主牌堆：获得遭遇牌堆第1张战利品牌。使用1张《感知》发动本牌时，转而翻开遭遇牌堆前3张战利品牌中选1张获得。"
$ws.Range("D6").Value = "使用1任意属性，或将1张手牌洗回主牌堆：选房间区1张陷阱牌移动到房间区任意位置。使用1张《感知》发动本牌时，可以转而选弃牌堆1张陷阱牌放到房间区任意位置。"
$ws.Range("D7").Value = "使用1任意属性，或将1张手牌洗回主牌堆：从遭遇牌堆翻开3张牌，获得其中的战利品牌。使用1张《敏捷》发动本牌时，可以额外翻开2张牌。"
$ws.Range("D8").Value = "消耗2时间：将主牌堆第1张怪物牌放在房间区任意非空列顶端，然后获得遭遇牌堆第1张战利品牌，再获得遗物牌堆顶的1张遗物牌。"
$ws.Range("D9").Value = "多选：①可重复 支付1金币为1张道具牌充1能。②支付3金币，从遗物牌堆翻开3张牌，选其中1张获得。③可重复 弃置1张任意战利品牌，获得1金币。"
$ws.Range("D10").Value = "消耗10时间，然后使用1任意属性或将1张手牌洗回主牌堆：升1级。使用1张《智力》发动本牌时，可以少消耗2时间。"

# Row heights for rows 3 & 5 grow because the new text wraps to one more line.
$ws.Rows.Item(3).RowHeight = 71.25
$ws.Rows.Item(5).RowHeight = 71.25

# View state: scroll back to the top (remove topLeftCell="A3") and move the
# active selection to D11.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D11").Select()

# Application window geometry.
$excel.Left = -120
$excel.Top = -120
$excel.Width = 29040
$excel.Height = 15720
